# Auto-generated edit script: updates cryptos price/volume table
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# These Price cells (column D) get numeric-looking new values; force
# text format first so Excel keeps them as text (matching original inlineStr type)
$textCells = @("D5","D6","D8","D9","D10","D14","D15","D16","D18","D20","D22","D23","D25","D26","D27","D30","D32","D37","D43","D48","D51")
foreach ($addr in $textCells) {
    $ws.Range($addr).NumberFormat = "@"
}

# Apply updated cell values
$ws.Range("D2").Value = "27.446.01"
$ws.Range("E2").Value = "  -2.16%  "
$ws.Range("D3").Value = "1.655.41"
$ws.Range("E3").Value = "  -1.95%  "
$ws.Range("E4").Value = "  -0.04%  "
$ws.Range("D5").Value = "213.80"
$ws.Range("E5").Value = "  -1.32%  "
$ws.Range("D6").Value = "0.512"
$ws.Range("E6").Value = "  -1.53%  "
$ws.Range("E7").Value = "  -0.04%  "
$ws.Range("D8").Value = "24.24"
$ws.Range("E8").Value = "  +0.56%  "
$ws.Range("D9").Value = "0.262"
$ws.Range("E9").Value = "  -1.07%  "
$ws.Range("D10").Value = "0.0616"
$ws.Range("E10").Value = "  -1.64%  "
$ws.Range("E11").Value = "  -0.63%  "
$ws.Range("D12").Value = "1.889.18"
$ws.Range("E12").Value = "  -1.98%  "
$ws.Range("D13").Value = "1.634.09"
$ws.Range("E13").Value = "  -3.24%  "
$ws.Range("D14").Value = "4.10"
$ws.Range("E14").Value = "  -2.08%  "
$ws.Range("D15").Value = "0.573"
$ws.Range("E15").Value = "  +2.47%  "
$ws.Range("D16").Value = "65.92"
$ws.Range("E16").Value = "  -1.46%  "
$ws.Range("D17").Value = "27.444.78"
$ws.Range("E17").Value = "  -2.01%  "
$ws.Range("D18").Value = "233.60"
$ws.Range("E18").Value = "  -6.78%  "
$ws.Range("D19").Value = "0.0₃0728"
$ws.Range("E19").Value = "  -1.95%  "
$ws.Range("D20").Value = "7.53"
$ws.Range("E20").Value = "  -2.13%  "
$ws.Range("E21").Value = "  -0.11%  "
$ws.Range("D22").Value = "4.40"
$ws.Range("E22").Value = "  -2.87%  "
$ws.Range("D23").Value = "9.33"
$ws.Range("E23").Value = "  -2.48%  "
$ws.Range("E24").Value = "  -1.37%  "
$ws.Range("D25").Value = "146.56"
$ws.Range("E25").Value = "  -0.72%  "
$ws.Range("D26").Value = "7.22"
$ws.Range("E26").Value = "  -1.85%  "
$ws.Range("D27").Value = "16.00"
$ws.Range("E27").Value = "  -3.03%  "
$ws.Range("E28").Value = "  +0.02%  "
$ws.Range("E29").Value = "  -2.11%  "
$ws.Range("D30").Value = "0.0498"
$ws.Range("E30").Value = "  -1.09%  "
$ws.Range("E31").Value = "  -4.23%  "
$ws.Range("D32").Value = "3.31"
$ws.Range("E32").Value = "  -2.33%  "
$ws.Range("D33").Value = "1.464.17"
$ws.Range("E33").Value = "  +2.61%  "
$ws.Range("E34").Value = "  -2.57%  "
$ws.Range("E35").Value = "  -3.90%  "
$ws.Range("E36").Value = "  -0.55%  "
$ws.Range("D37").Value = "0.913"
$ws.Range("E37").Value = "  -3.46%  "
$ws.Range("E38").Value = "  -3.10%  "
$ws.Range("E39").Value = "  -1.28%  "
$ws.Range("E40").Value = "  +0.08%  "
$ws.Range("E41").Value = "  -0.01%  "
$ws.Range("E42").Value = "  -0.89%  "
$ws.Range("D43").Value = "65.77"
$ws.Range("E43").Value = "  -5.46%  "
$ws.Range("E44").Value = "  -0.65%  "
$ws.Range("D45").Value = "1.797.51"
$ws.Range("E45").Value = "  -1.94%  "
$ws.Range("E46").Value = "  -1.84%  "
$ws.Range("E47").Value = "  +0.66%  "
$ws.Range("D48").Value = "88.51"
$ws.Range("E48").Value = "  -1.01%  "
$ws.Range("E49").Value = "  -4.02%  "
$ws.Range("E50").Value = "  -1.74%  "
$ws.Range("D51").Value = "7.82"
$ws.Range("E51").Value = "  -0.67%  "
